$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 250; existing rows 250-403 shift down to 251-404.
$ws.Rows("250:250").Insert()

# Populate the newly inserted row 250 with the new weekly record.
$ws.Range("A250").Value = 10
$ws.Range("B250").Value = "Vega Modelo de Temuco"
$ws.Range("C250").Value = "La Araucanía"
$ws.Range("D250").Value2 = 45126
$ws.Range("D250").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E250").Value = 9
$ws.Range("F250").Value = 100112039
$ws.Range("G250").Value = "Ciboulette"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 40
$ws.Range("K250").Value = 7000
$ws.Range("L250").Value = 7000
$ws.Range("M250").Value = 7000
$ws.Range("N250").Value = "`$/docena de atados"
$ws.Range("O250").Value = "Provincia de Cautín"
$ws.Range("P250").Value = 2333
$ws.Range("Q250").Value = 3
$ws.Range("R250").Value = "Hortaliza"
